$wb = $excel.ActiveWorkbook

# Update "zh-cn" sheet: Correspond Handoff Datetime (E) and Correspond Handback DateTime (H)
# Rows 2 and 3 (81ded9fc... and ad2cef13...) shared the same handoff/handback timestamps
# in the original report, so both rows must be refreshed together to keep them in sync.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-20 12:19:02"
$wsZh.Range("E3").Value = "2016-03-20 12:19:02"
$wsZh.Range("H2").Value = "2016-03-20 12:19:23"
$wsZh.Range("H3").Value = "2016-03-20 12:19:23"

# Update "de-de" sheet: Correspond Handoff Datetime (E) and Correspond Handback DateTime (H)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-20 12:19:06"
$wsDe.Range("E3").Value = "2016-03-20 12:19:06"
$wsDe.Range("H2").Value = "2016-03-20 12:19:28"
$wsDe.Range("H3").Value = "2016-03-20 12:19:28"
